$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(4,6,9,13,19,22,24,25,26,30,34,35,37,38,40,41,44,45,46,52,55,57,63,66,67,68,73,75,78,80,81,82,89,93,94,99,101,108,109,113,114,120,121,122,128,130,132,133,134,137,140,141,142,145,150,155,156,159,166,171,178,179,184,190)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
